# Testing slide (slide 4): add a new bullet/paragraph before "All Unit Tests passed."
# describing why unit tests were written after each function/method.
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(4)
$shape = $s.Shapes.Item(2)
$tr = $shape.TextFrame.TextRange

$newText = "Unit tests where created after each function/method before starting a new method/function. I chose to do this rather than do the test first because I had more experience doing it this way and to me its easier to write a test for something you have created than to write a test before that."

# Insert as a new paragraph immediately before the existing "All Unit Tests passed." paragraph
# (currently paragraph 2 of this placeholder).
$null = $tr.Paragraphs(2).InsertBefore($newText + "`r")
